$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.044463038444519
$ws.Range("B1").Value = 1.39484441280365
$ws.Range("C1").Value = 2.268747568130493
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.874749779701233
